$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-16 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-08-17 Saturday", 2) | Out-Null
$d.Content.Find.Execute("17÷8=2, 1", $true, $true, $false, $false, $false, $true, 1, $false, "92÷4=23, 0", 2) | Out-Null
$d.Content.Find.Execute("13÷9=1, 4", $true, $true, $false, $false, $false, $true, 1, $false, "53÷6=8, 5", 2) | Out-Null
$d.Content.Find.Execute("57÷8=7, 1", $true, $true, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2) | Out-Null
$d.Content.Find.Execute("69÷2=34, 1", $true, $true, $false, $false, $false, $true, 1, $false, "21÷5=4, 1", 2) | Out-Null
$d.Content.Find.Execute("11÷4=2, 3", $true, $true, $false, $false, $false, $true, 1, $false, "12÷4=3, 0", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $true, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $true, $false, $false, $false, $true, 1, $false, "99÷9=11, 0", 2) | Out-Null
$d.Content.Find.Execute("89÷9=9, 8", $true, $true, $false, $false, $false, $true, 1, $false, "82÷4=20, 2", 2) | Out-Null
$d.Content.Find.Execute("68÷2=34, 0", $true, $true, $false, $false, $false, $true, 1, $false, "70÷5=14, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷2=8, 1", $true, $true, $false, $false, $false, $true, 1, $false, "57÷3=19, 0", 2) | Out-Null
$d.Content.Find.Execute("98÷2=49, 0", $true, $true, $false, $false, $false, $true, 1, $false, "27÷6=4, 3", 2) | Out-Null
$d.Content.Find.Execute("20÷6=3, 2", $true, $true, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 2) | Out-Null
$d.Content.Find.Execute("64÷5=12, 4", $true, $true, $false, $false, $false, $true, 1, $false, "78÷4=19, 2", 2) | Out-Null
$d.Content.Find.Execute("50÷3=16, 2", $true, $true, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷6=3, 3", $true, $true, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("45÷8=5, 5", $true, $true, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 2) | Out-Null
$d.Content.Find.Execute("85÷2=42, 1", $true, $true, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2) | Out-Null
$d.Content.Find.Execute("10÷2=5, 0", $true, $true, $false, $false, $false, $true, 1, $false, "87÷9=9, 6", 2) | Out-Null
$d.Content.Find.Execute("70÷3=23, 1", $true, $true, $false, $false, $false, $true, 1, $false, "20÷2=10, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷4=13, 3", $true, $true, $false, $false, $false, $true, 1, $false, "30÷8=3, 6", 2) | Out-Null
$d.Content.Find.Execute("88÷9=9, 7", $true, $true, $false, $false, $false, $true, 1, $false, "15÷4=3, 3", 2) | Out-Null
$d.Content.Find.Execute("97÷7=13, 6", $true, $true, $false, $false, $false, $true, 1, $false, "68÷3=22, 2", 2) | Out-Null
$d.Content.Find.Execute("57÷7=8, 1", $true, $true, $false, $false, $false, $true, 1, $false, "64÷4=16, 0", 2) | Out-Null
$d.Content.Find.Execute("22÷5=4, 2", $true, $true, $false, $false, $false, $true, 1, $false, "77÷2=38, 1", 2) | Out-Null
$d.Content.Find.Execute("14÷8=1, 6", $true, $true, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2) | Out-Null
